# start route 110 trainers & route 103 eastern end trainers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "moves" header cells (column E) use the same distinct style as every
# other trainer-table header row already in the sheet (e.g. E456) - not a
# visual change (same non-bold Arial 10), just matches the existing pattern.
$movesHeaderStyle = $ws.Range("E456").Style

# --- route 110 section header ---
$ws.Range("A460").Value = "# route 110"

# TRAINER_KALEB
$ws.Range("A462").Value = "TRAINER_KALEB"
$ws.Range("A463").Value = "species"
$ws.Range("B463").Value = "lvl"
$ws.Range("C463").Value = "iv"
$ws.Range("D463").Value = "heldItem"
$ws.Range("E463").Value = "moves"
$ws.Range("E463").Style = $movesHeaderStyle
$ws.Range("F463").Value = "ability"
$ws.Range("G463").Value = "shiny"
$ws.Range("A464").Value = "Alolan_Vulpix"
$ws.Range("B464").Value = 27
$ws.Range("A465").Value = "Vulpix"
$ws.Range("B465").Value = 27
$ws.Range("G465").Value = 1

# TRAINER_ISABEL_1
$ws.Range("A467").Value = "TRAINER_ISABEL_1"
$ws.Range("A468").Value = "species"
$ws.Range("B468").Value = "lvl"
$ws.Range("C468").Value = "iv"
$ws.Range("D468").Value = "heldItem"
$ws.Range("E468").Value = "moves"
$ws.Range("E468").Style = $movesHeaderStyle
$ws.Range("F468").Value = "ability"
$ws.Range("G468").Value = "shiny"
$ws.Range("A469").Value = "Snubbull"
$ws.Range("B469").Value = 27
$ws.Range("A470").Value = "Clefairy"
$ws.Range("B470").Value = 28

# --- route 103 eastern end section header (replaces old "END" marker) ---
$ws.Range("A472").Value = "# route 103 eastern end"

# TRAINER_AMY_AND_LIV_1
$ws.Range("A474").Value = "TRAINER_AMY_AND_LIV_1"
$ws.Range("A475").Value = "species"
$ws.Range("B475").Value = "lvl"
$ws.Range("C475").Value = "iv"
$ws.Range("D475").Value = "heldItem"
$ws.Range("E475").Value = "moves"
$ws.Range("E475").Style = $movesHeaderStyle
$ws.Range("F475").Value = "ability"
$ws.Range("G475").Value = "shiny"
$ws.Range("A476").Value = "Plusle"
$ws.Range("B476").Value = 28
$ws.Range("A477").Value = "Minun"
$ws.Range("B477").Value = 28

# TRAINER_DAISY
$ws.Range("A479").Value = "TRAINER_DAISY"
$ws.Range("A480").Value = "species"
$ws.Range("B480").Value = "lvl"
$ws.Range("C480").Value = "iv"
$ws.Range("D480").Value = "heldItem"
$ws.Range("E480").Value = "moves"
$ws.Range("E480").Style = $movesHeaderStyle
$ws.Range("F480").Value = "ability"
$ws.Range("G480").Value = "shiny"
$ws.Range("A481").Value = "Budew"
$ws.Range("B481").Value = 27
$ws.Range("A482").Value = "Bayleef"
$ws.Range("B482").Value = 29

# restore final view state (scroll position / selection) similar to the source edit
$ws.Range("B483").Select()
